# admin password change done
#
# 1. "php artisan view:clear" paragraph gets an <w:rStyle w:val="pln"/>
#    added to its paragraph-mark run properties (w:pPr/w:rPr) - a side
#    effect of Word syncing the mark's rPr with the run that now precedes
#    a freshly-typed new paragraph.
# 2. A new paragraph "do this above 3 before backing up" is inserted right
#    after it (style "l2", no list numbering).
# 3. A new paragraph "do this when any changes is made in web.php" is
#    inserted after the first blank line following "php artisan optimize".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: locate the "php artisan view:clear" paragraph and rebuild it
# with the rStyle added to the paragraph mark's rPr. (InsertXML drops
# w:rStyle on *run* rPr, so we restore each run's character style with
# Range.Style right after.)
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "php artisan view:clear") {
        $target = $p
        break
    }
}

$r3 = $target.Range
$start3 = $r3.Start

$xml3 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:pPr>" +
          "<w:pStyle w:val='l2'/>" +
          "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='16'/></w:numPr>" +
          "<w:spacing w:before='0' w:after='0' w:afterAutospacing='0'/>" +
          "<w:rPr>" +
            "<w:rStyle w:val='pln'/>" +
            "<w:rFonts w:ascii='Consolas' w:hAnsi='Consolas'/>" +
            "<w:color w:val='2D2F31'/>" +
            "<w:sz w:val='18'/>" +
            "<w:szCs w:val='18'/>" +
          "</w:rPr>" +
        "</w:pPr>" +
        "<w:r><w:rPr><w:rFonts w:ascii='Consolas' w:hAnsi='Consolas'/><w:color w:val='2D2F31'/><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>php artisan view</w:t></w:r>" +
        "<w:r><w:rPr><w:rFonts w:ascii='Consolas' w:hAnsi='Consolas'/><w:color w:val='2D2F31'/><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>:</w:t></w:r>" +
        "<w:r><w:rPr><w:rFonts w:ascii='Consolas' w:hAnsi='Consolas'/><w:color w:val='2D2F31'/><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>clear</w:t></w:r>" +
        "</w:p>"
$r3.InsertXML($xml3)

# Re-apply the character styles that InsertXML stripped off the runs.
$seg1 = $d.Range($start3, $start3 + 16)           # "php artisan view"
$seg1.Style = "pln"
$seg2 = $d.Range($start3 + 16, $start3 + 17)      # ":"
$seg2.Style = "pun"
$seg3 = $d.Range($start3 + 17, $start3 + 22)      # "clear"
$seg3.Style = "pln"

# ---------------------------------------------------------------------
# Step 2: insert the new "do this above 3 before backing up" paragraph
# right after the paragraph we just rebuilt.
# ---------------------------------------------------------------------
$pClear = $d.Paragraphs.Item($target.Index)
$pClear.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item($target.Index + 1)
$xmlNew1 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
             "<w:pPr>" +
               "<w:pStyle w:val='l2'/>" +
               "<w:spacing w:before='0' w:after='0' w:afterAutospacing='0'/>" +
               "<w:rPr>" +
                 "<w:rFonts w:ascii='Consolas' w:eastAsiaTheme='minorEastAsia' w:hAnsi='Consolas' w:hint='eastAsia'/>" +
                 "<w:color w:val='2D2F31'/>" +
                 "<w:sz w:val='18'/>" +
                 "<w:szCs w:val='18'/>" +
               "</w:rPr>" +
             "</w:pPr>" +
             "<w:r>" +
               "<w:rPr>" +
                 "<w:rFonts w:ascii='Consolas' w:eastAsiaTheme='minorEastAsia' w:hAnsi='Consolas'/>" +
                 "<w:color w:val='2D2F31'/>" +
                 "<w:sz w:val='18'/>" +
                 "<w:szCs w:val='18'/>" +
               "</w:rPr>" +
               "<w:t>do this above 3 before backing up</w:t>" +
             "</w:r>" +
           "</w:p>"
$pNew1.Range.InsertXML($xmlNew1)

# ---------------------------------------------------------------------
# Step 3: insert "do this when any changes is made in web.php" right
# after the first blank paragraph that follows "php artisan optimize".
# ---------------------------------------------------------------------
$pOptimize = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "php artisan optimize") {
        $pOptimize = $p
        break
    }
}
$pBlank = $d.Paragraphs.Item($pOptimize.Index + 1)
$pBlank.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs.Item($pOptimize.Index + 2)
$xmlNew2 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
             "<w:pPr>" +
               "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/>" +
               "<w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr>" +
             "</w:pPr>" +
             "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>d</w:t></w:r>" +
             "<w:r><w:t>o this when any changes is made in web.php</w:t></w:r>" +
           "</w:p>"
$pNew2.Range.InsertXML($xmlNew2)

Write-Output "done"
